# "add test suite collection for login, register, and visit"
#
# Updates the sample/seed row on the "New Patient" sheet with a new test
# persona (Leonardo Di Caprio / relative Marshanda) and a renamed
# appointment type, bumps the min/max time-frame values, widens the
# AppointmentType column to fit the new text, and re-selects the last
# cell of the header row. Also nudges the "Random" helper sheet so its
# volatile RANDBETWEEN()-based day/year are regenerated.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("New Patient")

# --- Row 2 sample data -----------------------------------------------
# Set the shared-string-producing cells in the same order the authored
# workbook introduced them in, so the new entries land in the expected
# spots: AppointmentType, RelativesName, then GivenName/MiddleName/FamilyName.
$ws.Range("N2").Value = "General Medicine (New Patient)"
$ws.Range("M2").Value = "MARSHANDA"
$ws.Range("B2").Value = "LEONARDO"
$ws.Range("C2").Value = "DI"
$ws.Range("D2").Value = "CAPRIO"

# MinTimeFrameValue / MaxTimeFrameValue
$ws.Range("P2").Value = 4
$ws.Range("Q2").Value = 5

# --- Column width ------------------------------------------------------
# AppointmentType's column needs to widen to fit "General Medicine (New
# Patient)".
$ws.Columns.Item(14).ColumnWidth = 26.71

# --- Selection / view ---------------------------------------------------
# Scroll the window right and leave the last header cell selected.
$excel.ActiveWindow.ScrollColumn = 8
$ws.Range("Q2").Select() | Out-Null

# --- Random helper sheet -------------------------------------------------
# Force the volatile RANDBETWEEN()-driven Day/Year values to regenerate.
$ws2 = $wb.Worksheets.Item("Random")
$excel.CalculateFull() | Out-Null
